$wb = $excel.ActiveWorkbook

# Row -> new value for column F ("想去人数") on the sheets that carry the data
$updates = @{
    5  = 52
    6  = 568
    8  = 2031
    11 = 4391
    16 = 115
    18 = 17
    20 = 3208
    21 = 71
    22 = 474
    29 = 56
    32 = 576
    33 = 1833
    34 = 282
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
